{"js": "// \"Fixed Bug- Image sizes in GalleryActivity\" \u2014 rewrite the Gallery Activity\n// test-case write-up: merge the PASS/GalleryActivity run, mark test #2 (and\n// its FAIL notes) as struck-through/closed with a \"Fixed 6/25/2013\" note,\n// renumber the type test to #3, change its verdict text, add a brand-new\n// test #4 about loading images from a non-gallery source (Dropbox / Google\n// Drive), and relocate the _GoBack bookmark to the new final paragraph.\n//\n// The new body is expressed as exact OOXML (the same shape Word would round\n// trip through getOoxml()/insertOoxml()) and dropped in as a single\n// replace of the document body, which lets every paragraph/run/bookmark\n// move land exactly as specified instead of being reconstructed through a\n// long chain of ad-hoc insert/format calls.\n\nconst body = context.document.body;\n\nconst newBodyOoxml =\n  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Test Case Documentation:</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:t>Gallery Activity testing:</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>1.</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Test: Images of different sizes. Small size</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>PASS.  Image displays in GalleryActivity</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>2. Test: Image of different sizes, Large size.</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:strike/></w:rPr><w:t>FAIL. Image displays, but the button to accept is not on the screen. Need to scale the</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:strike/></w:rPr><w:t>Image down so that the button can also be displayed.</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Fixed 6/25/2013</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>NOTE: DURING THIS TEST IT WAS NOTED THAT IF THE DEVICE ROTATES, THE GALLERY ACTIVITY RETURNS TO THE ACTIVITY TO SELECT AN IMAGE. THIS SHOULD NOT HAPPEN.</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>' +\n  '<w:p><w:r><w:t>3</w:t></w:r><w:r><w:t>. Test: Image of different type.</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>PASS, LOADS BOTH .JPEG AND .PNG FROM GALLERY</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>' +\n  '<w:p><w:r><w:t>4. Load image from non-gallery source</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>FAIL. WILL NOT LOAD ANY IMAGE FROM DROPBOX</w:t></w:r></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>FAIL. WILL TRY TO LOAD ANY FILETYPE FROM GOOGLE DRIVE, BUT NOT SUCCEED</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>' +\n  '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>';\n\nconst flatOpc =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<?mso-application progid=\"Word.Document\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships></pkg:xmlData></pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + newBodyOoxml + '</w:body></w:document></pkg:xmlData></pkg:part>' +\n  '</pkg:package>';\n\nbody.insertOoxml(flatOpc, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# \"Fixed Bug- Image sizes in GalleryActivity\" \u2014 rewrite the Gallery Activity\n# test-case write-up: merge the PASS/GalleryActivity run, mark test #2 (and\n# its FAIL notes) as struck-through/closed with a \"Fixed 6/25/2013\" note,\n# renumber the type test to #3, change its verdict text, add a brand-new\n# test #4 about loading images from a non-gallery source (Dropbox / Google\n# Drive), and relocate the _GoBack bookmark to the new final paragraph.\n#\n# Each paragraph is built as a WordprocessingML <w:p> fragment (with the\n# w: namespace declared on every top-level fragment, since InsertXML parses\n# each root node independently) and the whole run of paragraphs is dropped\n# into the document in one shot via Range.InsertXML, which replaces that\n# range's contents \u2014 exactly mirroring the paragraph-by-paragraph shape of\n# the target document.\n\n$w = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n$p01 = \"<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Test Case Documentation:</w:t></w:r></w:p>\"\n$p02 = \"<w:p $w><w:pPr><w:rPr><w:u w:val=`\"single`\"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val=`\"single`\"/></w:rPr><w:t>Gallery Activity testing:</w:t></w:r></w:p>\"\n$p03 = \"<w:p $w><w:r><w:t>1.</w:t></w:r></w:p>\"\n$p04 = \"<w:p $w><w:r><w:t>Test: Images of different sizes. Small size</w:t></w:r></w:p>\"\n$p05 = \"<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>PASS.  Image displays in GalleryActivity</w:t></w:r></w:p>\"\n$p06 = \"<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>\"\n$p07 = \"<w:p $w><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>2. Test: Image of different sizes, Large size.</w:t></w:r></w:p>\"\n$p08 = \"<w:p $w><w:pPr><w:rPr><w:b/><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:strike/></w:rPr><w:t>FAIL. Image displays, but the button to accept is not on the screen. Need to scale the</w:t></w:r></w:p>\"\n$p09 = \"<w:p $w><w:pPr><w:rPr><w:b/><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:strike/></w:rPr><w:t>Image down so that the button can also be displayed.</w:t></w:r></w:p>\"\n$p10 = \"<w:p $w><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Fixed 6/25/2013</w:t></w:r></w:p>\"\n$p11 = \"<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>\"\n$p12 = \"<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>NOTE: DURING THIS TEST IT WAS NOTED THAT IF THE DEVICE ROTATES, THE GALLERY ACTIVITY RETURNS TO THE ACTIVITY TO SELECT AN IMAGE. THIS SHOULD NOT HAPPEN.</w:t></w:r></w:p>\"\n$p13 = \"<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>\"\n$p14 = \"<w:p $w><w:r><w:t>3</w:t></w:r><w:r><w:t>. Test: Image of different type.</w:t></w:r></w:p>\"\n$p15 = \"<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>PASS, LOADS BOTH .JPEG AND .PNG FROM GALLERY</w:t></w:r></w:p>\"\n$p16 = \"<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>\"\n$p17 = \"<w:p $w><w:r><w:t>4. Load image from non-gallery source</w:t></w:r></w:p>\"\n$p18 = \"<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>FAIL. WILL NOT LOAD ANY IMAGE FROM DROPBOX</w:t></w:r></w:p>\"\n$p19 = \"<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>FAIL. WILL TRY TO LOAD ANY FILETYPE FROM GOOGLE DRIVE, BUT NOT SUCCEED</w:t></w:r><w:bookmarkStart w:id=`\"0`\" w:name=`\"_GoBack`\"/><w:bookmarkEnd w:id=`\"0`\"/></w:p>\"\n$p20 = \"<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>\"\n\n$newXml = $p01 + $p02 + $p03 + $p04 + $p05 + $p06 + $p07 + $p08 + $p09 + $p10 + `\n          $p11 + $p12 + $p13 + $p14 + $p15 + $p16 + $p17 + $p18 + $p19 + $p20\n\n$d = $word.ActiveDocument\n$d.Content.InsertXML($newXml)\n"}
